# Updated CHE_grids model - 2025-08-21 08:22
#
# Re-applies the regenerated "existing_stock" values: the comm-out commodity,
# ncap_pasti (capacity factor), ncap_cost and description columns for a
# handful of hydro/solar rows were re-shuffled by the upstream VerveStacks
# build script. This reproduces the resulting cell values exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

# Row 10
$ws.Range("P10").Value = 'Aggregated Plant - IRENA Gap - relation/7933294-380_Missing Hydro Capacity'

# Row 11
$ws.Range("C11").Value = 'e_w240959264-220'
$ws.Range("E11").Value = 0.1117697096937668
$ws.Range("G11").Value = 3162.5
$ws.Range("P11").Value = 'Aggregated Plant - IRENA Gap - way/240959264-220_Missing Hydro Capacity'

# Row 12
$ws.Range("C12").Value = 'e_r7933294-380'
$ws.Range("E12").Value = 0.12488235719973945
$ws.Range("G12").Value = 3162.5000000000005
$ws.Range("P12").Value = 'Aggregated Plant - IRENA Gap - CH18-220_Missing Hydro Capacity'

# Row 13
$ws.Range("P13").Value = 'Aggregated Plant - IRENA Gap - way/238138373-380_Missing Hydro Capacity'

# Row 137
$ws.Range("P137").Value = 'Aggregated Plant - IRENA Gap - CHE_21_Missing Solar Capacity'

# Row 140
$ws.Range("P140").Value = 'Aggregated Plant - IRENA Gap - CHE_12_Missing Solar Capacity'

# Row 141
$ws.Range("P141").Value = 'Aggregated Plant - IRENA Gap - CHE_11_Missing Solar Capacity'

# Row 142
$ws.Range("P142").Value = 'Aggregated Plant - IRENA Gap - CHE_7_Missing Solar Capacity'

# Row 143
$ws.Range("P143").Value = 'Aggregated Plant - IRENA Gap - CHE_19_Missing Solar Capacity'

# Row 144
$ws.Range("P144").Value = 'Aggregated Plant - IRENA Gap - CHE_3_Missing Solar Capacity'

# Row 145
$ws.Range("P145").Value = 'Aggregated Plant - IRENA Gap - CHE_9_Missing Solar Capacity'

# Row 147
$ws.Range("P147").Value = 'Aggregated Plant - IRENA Gap - CHE_4_Missing Solar Capacity'

# Row 148
$ws.Range("P148").Value = 'Aggregated Plant - IRENA Gap - CHE_1_Missing Solar Capacity'

# Row 149
$ws.Range("P149").Value = 'Aggregated Plant - IRENA Gap - CHE_17_Missing Solar Capacity'

# Row 150
$ws.Range("C150").Value = 'elc_spv-CHE_0019'
$ws.Range("E150").Value = 0.16485344960649678
$ws.Range("P150").Value = 'Aggregated Plant - IRENA Gap - CHE_13_Missing Solar Capacity'

# Row 151
$ws.Range("C151").Value = 'elc_spv-CHE_0012'
$ws.Range("E151").Value = 0.13549669849969209
$ws.Range("P151").Value = 'Aggregated Plant - IRENA Gap - CHE_24_Missing Solar Capacity'

# Row 152
$ws.Range("C152").Value = 'elc_spv-CHE_0011'
$ws.Range("E152").Value = 0.16209575724687297
$ws.Range("P152").Value = 'Aggregated Plant - IRENA Gap - CHE_20_Missing Solar Capacity'

# Row 153
$ws.Range("C153").Value = 'elc_spv-CHE_0003'
$ws.Range("E153").Value = 0.16085025627375071
$ws.Range("P153").Value = 'Aggregated Plant - IRENA Gap - CHE_18_Missing Solar Capacity'

# Row 154
$ws.Range("C154").Value = 'elc_spv-CHE_0021'
$ws.Range("E154").Value = 0.15273795001145538
$ws.Range("P154").Value = 'Aggregated Plant - IRENA Gap - CHE_2_Missing Solar Capacity'

# Row 155
$ws.Range("C155").Value = 'elc_spv-CHE_0006'
$ws.Range("E155").Value = 0.21640319337561012
$ws.Range("P155").Value = 'Aggregated Plant - IRENA Gap - CHE_6_Missing Solar Capacity'

# Row 156
$ws.Range("C156").Value = 'elc_spv-CHE_0013'
$ws.Range("E156").Value = 0.17206733071733712
$ws.Range("P156").Value = 'Aggregated Plant - IRENA Gap - CHE_23_Missing Solar Capacity'

# Row 158
$ws.Range("C158").Value = 'elc_spv-CHE_0025'
$ws.Range("E158").Value = 0.13652468601509371
$ws.Range("P158").Value = 'Aggregated Plant - IRENA Gap - CHE_8_Missing Solar Capacity'

# Row 159
$ws.Range("C159").Value = 'elc_spv-CHE_0005'
$ws.Range("E159").Value = 0.20006982412215921
$ws.Range("P159").Value = 'Aggregated Plant - IRENA Gap - CHE_22_Missing Solar Capacity'

# Row 160
$ws.Range("P160").Value = 'Aggregated Plant - IRENA Gap - CHE_25_Missing Solar Capacity'

# Row 166
$ws.Range("C166").Value = 'elc_spv-CHE_0017'
$ws.Range("E166").Value = 0.15226887751132734

# Row 167
$ws.Range("C167").Value = 'elc_spv-CHE_0020'
$ws.Range("E167").Value = 0.15456128021356608

# Row 168
$ws.Range("C168").Value = 'elc_spv-CHE_0007'
$ws.Range("E168").Value = 0.16629376698088194

# Row 169
$ws.Range("C169").Value = 'elc_spv-CHE_0022'
$ws.Range("E169").Value = 0.21381383751804844

# Row 170
$ws.Range("C170").Value = 'elc_spv-CHE_0004'
$ws.Range("E170").Value = 0.19745398836539674

# Row 174
$ws.Range("C174").Value = 'elc_spv-CHE_0014'
$ws.Range("E174").Value = 0.18231505170803797
